{"js": "// Replace the three-digit x one-digit multiplication problems throughout\n// the document's table cells with a new set of problems, per the commit's\n// regenerated data. Every \"<old>=\" string is unique in the document, so a\n// direct search-and-replace for each pair is unambiguous.\nconst replacements = [\n  [\"556\u00d79=\", \"767\u00d78=\"],\n  [\"667\u00d73=\", \"590\u00d74=\"],\n  [\"325\u00d77=\", \"698\u00d73=\"],\n  [\"475\u00d73=\", \"181\u00d72=\"],\n  [\"462\u00d75=\", \"233\u00d76=\"],\n  [\"197\u00d79=\", \"482\u00d77=\"],\n  [\"260\u00d75=\", \"832\u00d74=\"],\n  [\"146\u00d74=\", \"134\u00d73=\"],\n  [\"681\u00d73=\", \"176\u00d77=\"],\n  [\"977\u00d73=\", \"215\u00d72=\"],\n  [\"564\u00d72=\", \"452\u00d74=\"],\n  [\"678\u00d72=\", \"403\u00d74=\"],\n  [\"495\u00d76=\", \"728\u00d78=\"],\n  [\"668\u00d78=\", \"856\u00d79=\"],\n  [\"479\u00d77=\", \"892\u00d73=\"],\n  [\"361\u00d75=\", \"715\u00d73=\"],\n  [\"207\u00d72=\", \"225\u00d76=\"],\n  [\"515\u00d78=\", \"489\u00d76=\"],\n  [\"186\u00d79=\", \"832\u00d77=\"],\n  [\"502\u00d72=\", \"260\u00d76=\"],\n  [\"640\u00d75=\", \"739\u00d72=\"],\n  [\"198\u00d74=\", \"287\u00d79=\"],\n  [\"901\u00d78=\", \"261\u00d74=\"],\n  [\"792\u00d75=\", \"475\u00d72=\"],\n  [\"220\u00d77=\", \"469\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems throughout\n# the document's table cells with a new set of problems, per the commit's\n# regenerated data. Every \"<old>=\" string is unique in the document, so a\n# direct Find/Replace for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"556\u00d79=\", \"767\u00d78=\"),\n    @(\"667\u00d73=\", \"590\u00d74=\"),\n    @(\"325\u00d77=\", \"698\u00d73=\"),\n    @(\"475\u00d73=\", \"181\u00d72=\"),\n    @(\"462\u00d75=\", \"233\u00d76=\"),\n    @(\"197\u00d79=\", \"482\u00d77=\"),\n    @(\"260\u00d75=\", \"832\u00d74=\"),\n    @(\"146\u00d74=\", \"134\u00d73=\"),\n    @(\"681\u00d73=\", \"176\u00d77=\"),\n    @(\"977\u00d73=\", \"215\u00d72=\"),\n    @(\"564\u00d72=\", \"452\u00d74=\"),\n    @(\"678\u00d72=\", \"403\u00d74=\"),\n    @(\"495\u00d76=\", \"728\u00d78=\"),\n    @(\"668\u00d78=\", \"856\u00d79=\"),\n    @(\"479\u00d77=\", \"892\u00d73=\"),\n    @(\"361\u00d75=\", \"715\u00d73=\"),\n    @(\"207\u00d72=\", \"225\u00d76=\"),\n    @(\"515\u00d78=\", \"489\u00d76=\"),\n    @(\"186\u00d79=\", \"832\u00d77=\"),\n    @(\"502\u00d72=\", \"260\u00d76=\"),\n    @(\"640\u00d75=\", \"739\u00d72=\"),\n    @(\"198\u00d74=\", \"287\u00d79=\"),\n    @(\"901\u00d78=\", \"261\u00d74=\"),\n    @(\"792\u00d75=\", \"475\u00d72=\"),\n    @(\"220\u00d77=\", \"469\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
